# Auto-generated edit script
# Updates rows 622-681 (columns A-T) of the active worksheet to reflect the
# newly published weekly price data (week of 2021-09-22). The previously
# existing weekly blocks shift down by four rows, and the oldest block
# (2021-05-07) is duplicated as new trailing rows 678-681, extending the
# used range from A1:T677 to A1:T681.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 622

$data = @(
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44461, 13, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Especial', 40, 17000, 17000, 17000, '$/caja 10 unidades', 'Ecuador', 1700, 10),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44461, 13, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Primera', 45, 17000, 17000, 17000, '$/caja 12 unidades', 'Ecuador', 1417, 12),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44461, 13, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Segunda', 30, 17000, 17000, 17000, '$/caja 14 unidades', 'Ecuador', 1214, 14),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44461, 13, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Tercera', 25, 17000, 17000, 17000, '$/caja 16 unidades', 'Ecuador', 1062, 16),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44357, 13, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Especial', 108, 13500, 14500, 14000, '$/caja 10 unidades', 'Ecuador', 1400, 10),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44357, 13, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Primera', 216, 13500, 14500, 14000, '$/caja 12 unidades', 'Ecuador', 1167, 12),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44357, 13, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Segunda', 216, 13500, 14500, 14000, '$/caja 14 unidades', 'Ecuador', 1000, 14),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44203, 13, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Especial', 50, 14000, 14000, 14000, '$/caja 10 unidades', 'Ecuador', 1400, 10),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44203, 13, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Primera', 50, 14000, 14000, 14000, '$/caja 12 unidades', 'Ecuador', 1167, 12),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44203, 13, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Segunda', 35, 14000, 14000, 14000, '$/caja 14 unidades', 'Ecuador', 1000, 14),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44203, 13, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Tercera', 40, 14000, 14000, 14000, '$/caja 16 unidades', 'Ecuador', 875, 16),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44162, 13, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Especial', 25, 22000, 22000, 22000, '$/caja 10 unidades', 'Ecuador', 2200, 10),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44162, 13, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Primera', 40, 22000, 22000, 22000, '$/caja 12 unidades', 'Ecuador', 1833, 12),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44162, 13, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Segunda', 35, 22000, 22000, 22000, '$/caja 14 unidades', 'Ecuador', 1571, 14),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44162, 13, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Tercera', 40, 22000, 22000, 22000, '$/caja 16 unidades', 'Ecuador', 1375, 16),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44410, 13, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Especial', 35, 17000, 18000, 17571, '$/caja 10 unidades', 'Ecuador', 1757, 10),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44410, 13, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Primera', 35, 17000, 18000, 17571, '$/caja 12 unidades', 'Ecuador', 1464, 12),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44410, 13, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Segunda', 25, 17000, 18000, 17600, '$/caja 14 unidades', 'Ecuador', 1257, 14),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44410, 13, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Tercera', 35, 17000, 18000, 17571, '$/caja 16 unidades', 'Ecuador', 1098, 16),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44411, 13, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Especial', 25, 16000, 16000, 16000, '$/caja 10 unidades', 'Ecuador', 1600, 10),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44411, 13, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Primera', 30, 16000, 16000, 16000, '$/caja 12 unidades', 'Ecuador', 1333, 12),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44411, 13, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Segunda', 30, 16000, 16000, 16000, '$/caja 14 unidades', 'Ecuador', 1143, 14),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44411, 13, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Tercera', 25, 16000, 16000, 16000, '$/caja 16 unidades', 'Ecuador', 1000, 16),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44257, 13, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Especial', 45, 14000, 14000, 14000, '$/caja 10 unidades', 'Ecuador', 1400, 10),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44257, 13, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Primera', 60, 14000, 14000, 14000, '$/caja 12 unidades', 'Ecuador', 1167, 12),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44257, 13, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Segunda', 50, 14000, 14000, 14000, '$/caja 14 unidades', 'Ecuador', 1000, 14),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44257, 13, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Tercera', 55, 14000, 14000, 14000, '$/caja 16 unidades', 'Ecuador', 875, 16),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44244, 13, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Especial', 40, 14000, 14000, 14000, '$/caja 10 unidades', 'Ecuador', 1400, 10),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44244, 13, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Primera', 45, 14000, 14000, 14000, '$/caja 12 unidades', 'Ecuador', 1167, 12),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44244, 13, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Segunda', 50, 14000, 14000, 14000, '$/caja 14 unidades', 'Ecuador', 1000, 14),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44244, 13, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Tercera', 40, 14000, 14000, 14000, '$/caja 16 unidades', 'Ecuador', 875, 16),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44176, 13, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Primera', 30, 17000, 17000, 17000, '$/caja 12 unidades', 'Ecuador', 1417, 12),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44176, 13, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Segunda', 40, 17000, 17000, 17000, '$/caja 14 unidades', 'Ecuador', 1214, 14),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44176, 13, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Tercera', 35, 17000, 17000, 17000, '$/caja 16 unidades', 'Ecuador', 1062, 16),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44239, 13, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Especial', 108, 14000, 15000, 14500, '$/caja 10 unidades', 'Ecuador', 1450, 10),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44239, 13, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Primera', 216, 14000, 15000, 14500, '$/caja 12 unidades', 'Ecuador', 1208, 12),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44239, 13, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Segunda', 216, 14000, 15000, 14500, '$/caja 14 unidades', 'Ecuador', 1036, 14),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44376, 13, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Especial', 40, 16000, 17000, 16625, '$/caja 10 unidades', 'Ecuador', 1662, 10),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44376, 13, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Primera', 35, 16000, 17000, 16429, '$/caja 12 unidades', 'Ecuador', 1369, 12),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44376, 13, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Segunda', 30, 16000, 17000, 16333, '$/caja 14 unidades', 'Ecuador', 1167, 14),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44376, 13, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Tercera', 40, 16000, 17000, 16375, '$/caja 16 unidades', 'Ecuador', 1023, 16),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44292, 13, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Especial', 55, 14000, 15000, 14545, '$/caja 10 unidades', 'Ecuador', 1454, 10),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44292, 13, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Primera', 55, 14000, 15000, 14545, '$/caja 12 unidades', 'Ecuador', 1212, 12),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44292, 13, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Segunda', 45, 14000, 15000, 14444, '$/caja 14 unidades', 'Ecuador', 1032, 14),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44292, 13, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Tercera', 45, 14000, 15000, 14444, '$/caja 16 unidades', 'Ecuador', 903, 16),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44358, 13, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Especial', 108, 13500, 14500, 14000, '$/caja 10 unidades', 'Ecuador', 1400, 10),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44358, 13, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Primera', 108, 13500, 14500, 14000, '$/caja 12 unidades', 'Ecuador', 1167, 12),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44358, 13, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Segunda', 150, 7500, 8000, 7750, '$/caja 7 unidades', 'Ecuador', 1107, 7),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44211, 13, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Especial', 20, 13000, 13000, 13000, '$/caja 10 unidades', 'Ecuador', 1300, 10),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44211, 13, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Primera', 30, 13000, 13000, 13000, '$/caja 12 unidades', 'Ecuador', 1083, 12),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44211, 13, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Segunda', 25, 13000, 13000, 13000, '$/caja 14 unidades', 'Ecuador', 929, 14),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44211, 13, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Tercera', 25, 13000, 13000, 13000, '$/caja 16 unidades', 'Ecuador', 812, 16),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44425, 13, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Especial', 45, 17000, 18000, 17444, '$/caja 10 unidades', 'Ecuador', 1744, 10),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44425, 13, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Primera', 55, 17000, 18000, 17455, '$/caja 12 unidades', 'Ecuador', 1455, 12),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44425, 13, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Segunda', 35, 17000, 18000, 17571, '$/caja 14 unidades', 'Ecuador', 1255, 14),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44425, 13, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Tercera', 30, 17000, 18000, 17667, '$/caja 16 unidades', 'Ecuador', 1104, 16),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44323, 13, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Especial', 25, 15000, 15000, 15000, '$/caja 10 unidades', 'Ecuador', 1500, 10),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44323, 13, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Primera', 30, 15000, 15000, 15000, '$/caja 12 unidades', 'Ecuador', 1250, 12),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44323, 13, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Segunda', 25, 15000, 15000, 15000, '$/caja 14 unidades', 'Ecuador', 1071, 14),
    @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44323, 13, 'Fruta', 100108, 'Tropicales y subtropicales', 100108005, 'Piña', 'Caramelo', 'Tercera', 30, 15000, 15000, 15000, '$/caja 16 unidades', 'Ecuador', 938, 16)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $r = $startRow + $i
    for ($j = 0; $j -lt $row.Length; $j++) {
        $ws.Cells.Item($r, $j + 1).Value = $row[$j]
    }
}

# Column D ("Fecha") uses a date number format; make sure the newly
# written cells (including the freshly appended rows 678-681) keep the
# same date formatting as the rest of the column.
for ($r = $startRow; $r -le ($startRow + $data.Length - 1); $r++) {
    $ws.Cells.Item($r, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
}
